$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Flip the "No" Runmode flags in column E (rows 2-29) over to "Yes"
#    so every prod sanity testcase is enabled for execution.
# ------------------------------------------------------------------
for ($r = 2; $r -le 29; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    if ($cell.Value() -eq "No") {
        $cell.Value = "Yes"
    }
}

# ------------------------------------------------------------------
# 2. Widen every conditional-format rule on column E that only covered
#    E2 or E2:E6 so it covers the full E2:E29 data range.
# ------------------------------------------------------------------
$target = $ws.Range("E2:E29")
$fcs = $ws.Range("E2").FormatConditions
$count = $fcs.Count()
for ($i = 1; $i -le $count; $i++) {
    $fc = $fcs.Item($i)
    $addr = $fc.AppliesTo().Address()
    if ($addr -eq "`$E`$2" -or $addr -eq "`$E`$2:`$E`$6") {
        $fc.ModifyAppliesToRange($target)
    }
}

# ------------------------------------------------------------------
# 3. Move the view/selection to match: scrolled up a bit, E2:E29 selected.
# ------------------------------------------------------------------
$excel.Goto($ws.Range("A7"), $false)
$target.Select()
